$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Find the last populated row in column A (xlUp = -4162) and append the new
# draw result directly beneath it, the way the daily auto-update job does.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
# Every column in this sheet stores plain text (dates, phase codes, and
# results all look numeric but must stay text, matching the rest of the log).
$rowRange.NumberFormat = "@"

$ws.Cells($newRow, 1).Value = "2025-10-23"
$ws.Cells($newRow, 2).Value = "Pick 4"
$ws.Cells($newRow, 3).Value = "251023"
$ws.Cells($newRow, 4).Value = "3-6-3-6"
$ws.Cells($newRow, 5).Value = "2025-10-23T21:38:22.853+04:00"
